# chore: update Sheets via scheduled runner
# Refresh currentAveragePrice / Leve profit figures (columns H-N) for a set
# of leve rows across the ALC, ARM, BSM, CRP, CUL, LTW and WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 70
$ws.Range("H70").Value = 2838.4707
$ws.Range("I70").Value = 1999
$ws.Range("J70").Value = 2890.9375
$ws.Range("K70").Value = 5997
$ws.Range("L70").Value = 8672.8125
$ws.Range("M70").Value = -5727
$ws.Range("N70").Value = -9212.8125

# Row 73
$ws.Range("H73").Value = 2838.4707
$ws.Range("I73").Value = 1999
$ws.Range("J73").Value = 2890.9375
$ws.Range("K73").Value = 5997
$ws.Range("L73").Value = 8672.8125
$ws.Range("M73").Value = -5061
$ws.Range("N73").Value = -10544.8125

# Row 92
$ws.Range("H92").Value = 145.58333
$ws.Range("I92").Value = 99.7
$ws.Range("J92").Value = 375
$ws.Range("K92").Value = 99.7
$ws.Range("L92").Value = 375
$ws.Range("M92").Value = 1148.3
$ws.Range("N92").Value = -2871

# Row 112
$ws.Range("H112").Value = 2500
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

# Row 138
$ws.Range("H138").Value = 3329.5
$ws.Range("J138").Value = 3313
$ws.Range("L138").Value = 9939
$ws.Range("N138").Value = -20219

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 1645.2593
$ws.Range("I32").Value = 757.7619
$ws.Range("J32").Value = 4751.5
$ws.Range("K32").Value = 757.7619
$ws.Range("L32").Value = 4751.5
$ws.Range("M32").Value = -470.7619
$ws.Range("N32").Value = -5325.5

# Row 95
$ws.Range("H95").Value = 47552
$ws.Range("J95").Value = 47552
$ws.Range("L95").Value = 47552
$ws.Range("N95").Value = -53044

# Row 97
$ws.Range("H97").Value = 1089.5
$ws.Range("I97").Value = 1089.5
$ws.Range("K97").Value = 1089.5
$ws.Range("M97").Value = -593.5

# Row 101
$ws.Range("H101").Value = 35000
$ws.Range("J101").Value = 35000
$ws.Range("L101").Value = 35000
$ws.Range("N101").Value = -41490

# Row 132
$ws.Range("H132").Value = 4982.4116
$ws.Range("I132").Value = 4793.8125
$ws.Range("K132").Value = 14381.4375
$ws.Range("M132").Value = -11851.4375

$ws = $wb.Worksheets.Item("BSM")
# Row 29
$ws.Range("H29").Value = 418.5
$ws.Range("I29").Value = 418.5
$ws.Range("K29").Value = 418.5
$ws.Range("M29").Value = -129.5

# Row 36
$ws.Range("H36").Value = 2215.2
$ws.Range("I36").Value = 2215.2
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 2215.2
$ws.Range("L36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("N36").Value = -1681.2

# Row 58
$ws.Range("H58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("L58").ClearContents()
$ws.Range("N58").Value = 0

# Row 59
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").ClearContents()
$ws.Range("N59").Value = 0

# Row 60
$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").ClearContents()
$ws.Range("N60").Value = 0

# Row 75
$ws.Range("H75").Value = 65000
$ws.Range("J75").Value = 65000
$ws.Range("L75").Value = 65000
$ws.Range("N75").Value = -66872

# Row 78
$ws.Range("H78").Value = 65000
$ws.Range("J78").Value = 65000
$ws.Range("L78").Value = 195000
$ws.Range("N78").Value = -204360

# Row 105
$ws.Range("H105").Value = 2416.4285
$ws.Range("I105").Value = 2152.5
$ws.Range("K105").Value = 2152.5
$ws.Range("M105").Value = -405.5

# Row 134
$ws.Range("H134").Value = 1938.8
$ws.Range("I134").Value = 1938.8
$ws.Range("K134").Value = 5816.4
$ws.Range("M134").Value = -3281.4

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 4045.4285
$ws.Range("J7").Value = 710.4
$ws.Range("L7").Value = 710.4
$ws.Range("N7").Value = -936.4

# Row 16
$ws.Range("H16").Value = 510.7857
$ws.Range("J16").Value = 665
$ws.Range("L16").Value = 665
$ws.Range("N16").Value = -1239

# Row 52
$ws.Range("H52").Value = 60000
$ws.Range("J52").Value = 60000
$ws.Range("L52").Value = 60000
$ws.Range("N52").Value = -60588

# Row 58
$ws.Range("H58").Value = 2816.2778
$ws.Range("I58").Value = 2262.1428
$ws.Range("K58").Value = 2262.1428
$ws.Range("M58").Value = -2059.1428

# Row 113
$ws.Range("H113").Value = 510.7857
$ws.Range("J113").Value = 665
$ws.Range("L113").Value = 665
$ws.Range("N113").Value = -5005

# Row 136
$ws.Range("H136").Value = 2816.2778
$ws.Range("I136").Value = 2262.1428
$ws.Range("K136").Value = 6786.428400000001
$ws.Range("M136").Value = -4236.428400000001

# Row 137
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").ClearContents()
$ws.Range("N137").Value = 0

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 74.43478399999999
$ws.Range("I2").Value = 161.5
$ws.Range("J2").Value = 28
$ws.Range("K2").Value = 969
$ws.Range("L2").Value = 168
$ws.Range("M2").Value = -856
$ws.Range("N2").Value = -394

# Row 4
$ws.Range("H4").Value = 100099.75

# Row 5
$ws.Range("H5").Value = 896.1429000000001
$ws.Range("I5").Value = 970.9286
$ws.Range("J5").Value = 846.2857
$ws.Range("K5").Value = 2912.7858
$ws.Range("L5").Value = 2538.8571
$ws.Range("M5").Value = -2800.7858
$ws.Range("N5").Value = -2762.8571

# Row 7
$ws.Range("H7").Value = 413.8889
$ws.Range("I7").Value = 285
$ws.Range("J7").Value = 575
$ws.Range("K7").Value = 855
$ws.Range("L7").Value = 1725
$ws.Range("M7").Value = -743
$ws.Range("N7").Value = -1949

# Row 74
$ws.Range("H74").Value = 8005
$ws.Range("I74").Value = 4995
$ws.Range("J74").Value = 8757.5
$ws.Range("K74").Value = 14985
$ws.Range("L74").Value = 26272.5
$ws.Range("M74").Value = -13924
$ws.Range("N74").Value = -28394.5

# Row 77
$ws.Range("H77").Value = 8005
$ws.Range("I77").Value = 4995
$ws.Range("J77").Value = 8757.5
$ws.Range("K77").Value = 44955
$ws.Range("L77").Value = 78817.5
$ws.Range("M77").Value = -39651
$ws.Range("N77").Value = -89425.5

# Row 107
$ws.Range("H107").Value = 301
$ws.Range("I107").Value = 301
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 903
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = 1017

# Row 121
$ws.Range("H121").Value = 199.66667

# Row 132
$ws.Range("H132").Value = 2507.4614
$ws.Range("J132").Value = 3299.8
$ws.Range("L132").Value = 29698.2
$ws.Range("N132").Value = -34758.2

# Row 135
$ws.Range("H135").Value = 896.1429000000001
$ws.Range("I135").Value = 970.9286
$ws.Range("J135").Value = 846.2857
$ws.Range("K135").Value = 8738.357399999999
$ws.Range("L135").Value = 7616.571300000001
$ws.Range("M135").Value = -6203.357399999999
$ws.Range("N135").Value = -12686.5713

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 6799.25
$ws.Range("I46").Value = 4842
$ws.Range("K46").Value = 4842
$ws.Range("M46").Value = -4654

$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Range("H136").Value = 3216.4348
$ws.Range("I136").Value = 2704.5881
$ws.Range("K136").Value = 8113.7643
$ws.Range("M136").Value = -5563.7643

Write-Host "Updated profit figures on ALC, ARM, BSM, CRP, CUL, LTW, WVR"
